# Work on the workbook that is already open
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Col box mod")

# --- Sheet "Tabelle1": zoom 90 -> 120, keep selection on D50 ---
[void]$ws1.Activate()
[void]$ws1.Range("D50").Select()
$excel.ActiveWindow.Zoom = 120

# --- Sheet "Col box mod": zoom 90 -> 120, selection moves to D22 ---
[void]$ws2.Activate()
[void]$ws2.Range("D22").Select()
$excel.ActiveWindow.Zoom = 120

# New narrow column G (rendered width ~6.24 "characters")
$ws2.Range("G1").ColumnWidth = 5.406666666666667

# New data block in rows 26-29 (columns B, C, D)
$ws2.Range("B26").Value = 1
$ws2.Range("C26").Value = "Hall"
$ws2.Range("D26").Value = "1000, 2000"

$ws2.Range("B27").Value = 2
$ws2.Range("C27").Value = "Hall"
$ws2.Range("D27").Value = 3000
$ws2.Rows.Item(27).RowHeight = 15

$ws2.Range("B28").Value = 3
$ws2.Range("C28").Value = "Long"
$ws2.Range("D28").Value = "1000, 2000"

$ws2.Range("B29").Value = 4
$ws2.Range("C29").Value = "Long"
$ws2.Range("D29").Value = 3000

# Leave the originally-selected tab active, as in the source file
[void]$ws2.Activate()
[void]$ws2.Range("D22").Select()

Write-Host "edit applied"
